# account_bank_statement_import_adyen test file migration
# ---------------------------------------------------------------------------
# The Adyen "credit fees" test fixture is re-denominated from EUR to USD
# (the Gross/Net Currency columns, K and O) and the sheet's on-screen
# view state (the active cell / scroll position) is updated to reflect
# where the author was working when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Currency re-denomination: every "EUR" value becomes "USD".
#    (The "Gross Currency"/"Net Currency" columns hold a handful of GBP
#    rows too -- those are left untouched, only literal "EUR" cells move.)
# ---------------------------------------------------------------------
$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -eq "EUR") {
            $cell.Value = "USD"
        }
    }
}

# ---------------------------------------------------------------------
# 2) View state: the active cell moves to P38 and the window scrolls so
#    row 10 is pinned at the top-left.
# ---------------------------------------------------------------------
$ws.Range("P38").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1

# Tab-area/horizontal-scrollbar split ratio on the workbook window.
$excel.ActiveWindow.TabRatio = 0.5
